# Updated cryptos list (Price / Volume(1h) columns) to match the
# latest scrape — mirrors the GitHub Actions commit that refreshed
# the coinranking.com snapshot on 2024-09-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.794.89'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').Value = '2.341.16'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''539.48'
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('D6').Value = '''134.02'
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('E7').Value = '  +0.73%  '
$ws.Range('E8').Value = '  +6.19%  '
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('E10').Value = '  +1.85%  '
$ws.Range('E11').Value = '  -1.64%  '
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('D14').Value = '2.757.52'
$ws.Range('E14').Value = '  +0.22%  '
$ws.Range('D15').Value = '57.743.85'
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = '2.332.79'
$ws.Range('E17').Value = '  -1.17%  '
$ws.Range('D18').Value = '''10.67'
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('D19').Value = '''4.29'
$ws.Range('E19').Value = '  +2.06%  '
$ws.Range('D20').Value = '''328.31'
$ws.Range('E20').Value = '  -2.11%  '
$ws.Range('E21').Value = '  -1.27%  '
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').Value = '''62.76'
$ws.Range('E23').Value = '  +1.60%  '
$ws.Range('D24').Value = '''0.163'
$ws.Range('E24').Value = '  -3.80%  '
$ws.Range('D25').Value = '''0.998'
$ws.Range('E25').Value = '  +0.39%  '
$ws.Range('D26').Value = '''8.30'
$ws.Range('E26').Value = '  -1.87%  '
$ws.Range('E27').Value = '  -6.03%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').Value = '''169.91'
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('E30').Value = '  -0.40%  '
$ws.Range('D31').Value = '''6.11'
$ws.Range('E31').Value = '  -0.94%  '
$ws.Range('E32').Value = '  -0.30%  '
$ws.Range('E33').Value = '  -1.45%  '
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  +0.87%  '
$ws.Range('E36').Value = '  +1.14%  '
$ws.Range('E37').Value = '  -2.14%  '
$ws.Range('D38').Value = '''1.59'
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('D39').Value = '''39.07'
$ws.Range('E39').Value = '  -0.73%  '
$ws.Range('D40').Value = '''141.51'
$ws.Range('E40').Value = '  -5.08%  '
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('D42').Value = '''3.63'
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('D43').Value = '''286.89'
$ws.Range('E43').Value = '  +1.67%  '
$ws.Range('E44').Value = '  +1.69%  '
$ws.Range('E45').Value = '  +0.79%  '
$ws.Range('D46').Value = '''19.13'
$ws.Range('E46').Value = '  -1.23%  '
$ws.Range('E47').Value = '  +0.80%  '
$ws.Range('E48').Value = '  +1.17%  '
$ws.Range('D49').Value = '''0.378'
$ws.Range('E49').Value = '  -1.06%  '
$ws.Range('D50').Value = '''11.07'
$ws.Range('E50').Value = '  +0.52%  '
$ws.Range('E51').Value = '  +0.99%  '